# Auto-generated edit script: updates Leve profit figures across 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 37444
$ws.Range("J68").Value = 37444
$ws.Range("L68").Value = 37444
$ws.Range("N68").Value = -38942
$ws.Range("H71").Value = 37444
$ws.Range("J71").Value = 37444
$ws.Range("L71").Value = 112332
$ws.Range("N71").Value = -119820
$ws.Range("H107").Value = 36106.32
$ws.Range("I107").Value = 38768.117
$ws.Range("J107").Value = 1503
$ws.Range("K107").Value = 38768.117
$ws.Range("L107").Value = 1503
$ws.Range("M107").Value = -36848.117
$ws.Range("N107").Value = -5343
$ws.Range("H132").Value = 2895.9167
$ws.Range("I132").Value = 1411.0667
$ws.Range("J132").Value = 5370.6665
$ws.Range("K132").Value = 4233.2001
$ws.Range("L132").Value = 16111.9995
$ws.Range("M132").Value = -1703.2001
$ws.Range("N132").Value = -21171.9995
$ws.Range("H134").Value = 38633.8
$ws.Range("J134").Value = 38633.8
$ws.Range("L134").Value = 38633.8
$ws.Range("N134").Value = -48773.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22468.97
$ws.Range("I32").Value = 4617.755
$ws.Range("J32").Value = 66204.45
$ws.Range("K32").Value = 4617.755
$ws.Range("L32").Value = 66204.45
$ws.Range("M32").Value = -4330.755
$ws.Range("N32").Value = -66778.45
$ws.Range("H74").Value = 774.93335
$ws.Range("I74").Value = 742
$ws.Range("J74").Value = 840.8
$ws.Range("K74").Value = 742
$ws.Range("L74").Value = 840.8
$ws.Range("M74").Value = 132
$ws.Range("N74").Value = -2588.8
$ws.Range("H77").Value = 774.93335
$ws.Range("I77").Value = 742
$ws.Range("J77").Value = 840.8
$ws.Range("K77").Value = 3710
$ws.Range("L77").Value = 4204
$ws.Range("M77").Value = 658
$ws.Range("N77").Value = -12940
$ws.Range("H97").Value = 2171.805
$ws.Range("I97").Value = 1560.091
$ws.Range("J97").Value = 2880.1052
$ws.Range("K97").Value = 1560.091
$ws.Range("L97").Value = 2880.1052
$ws.Range("M97").Value = -1064.091
$ws.Range("N97").Value = -3872.1052
$ws.Range("H110").Value = 2849.1667
$ws.Range("I110").Value = 3408.1667
$ws.Range("J110").Value = 1172.1666
$ws.Range("K110").Value = 3408.1667
$ws.Range("L110").Value = 1172.1666
$ws.Range("M110").Value = -1363.1667
$ws.Range("N110").Value = -5262.1666
$ws.Range("H122").Value = 1878.8148
$ws.Range("I122").Value = 1621.75
$ws.Range("J122").Value = 2613.2856
$ws.Range("K122").Value = 4865.25
$ws.Range("L122").Value = 7839.8568
$ws.Range("M122").Value = -2415.25
$ws.Range("N122").Value = -12739.8568
$ws.Range("H132").Value = 1921.2727
$ws.Range("I132").Value = 1892.5312
$ws.Range("K132").Value = 5677.5936
$ws.Range("M132").Value = -3147.5936
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1842.9131
$ws.Range("I86").Value = 2063
$ws.Range("J86").Value = 1339.8572
$ws.Range("K86").Value = 2063
$ws.Range("L86").Value = 1339.8572
$ws.Range("M86").Value = -940
$ws.Range("N86").Value = -3585.8572
$ws.Range("H89").Value = 1842.9131
$ws.Range("I89").Value = 2063
$ws.Range("J89").Value = 1339.8572
$ws.Range("K89").Value = 10315
$ws.Range("L89").Value = 6699.286
$ws.Range("M89").Value = -4699
$ws.Range("N89").Value = -17931.286
$ws.Range("H94").Value = 1168.5
$ws.Range("I94").Value = 1041.1666
$ws.Range("J94").Value = 1455
$ws.Range("K94").Value = 1041.1666
$ws.Range("L94").Value = 1455
$ws.Range("M94").Value = -590.1666
$ws.Range("N94").Value = -2357
$ws.Range("H99").Value = 2437.1428
$ws.Range("I99").Value = 1592.1111
$ws.Range("J99").Value = 3070.9167
$ws.Range("K99").Value = 1592.1111
$ws.Range("L99").Value = 3070.9167
$ws.Range("M99").Value = -94.11110000000008
$ws.Range("N99").Value = -6066.9167
$ws.Range("H134").Value = 1831.6923
$ws.Range("I134").Value = 1978.6666
$ws.Range("J134").Value = 1214.4
$ws.Range("K134").Value = 5935.9998
$ws.Range("L134").Value = 3643.2
$ws.Range("M134").Value = -3400.9998
$ws.Range("N134").Value = -8713.200000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1630.8387
$ws.Range("I58").Value = 1553.5
$ws.Range("K58").Value = 1553.5
$ws.Range("M58").Value = -1350.5
$ws.Range("H136").Value = 1630.8387
$ws.Range("I136").Value = 1553.5
$ws.Range("K136").Value = 4660.5
$ws.Range("M136").Value = -2110.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3261.4243
$ws.Range("I122").Value = 240.5
$ws.Range("J122").Value = 3456.3225
$ws.Range("K122").Value = 2164.5
$ws.Range("L122").Value = 31106.9025
$ws.Range("M122").Value = 285.5
$ws.Range("N122").Value = -36006.9025
$ws.Range("H126").Value = 1748.8235
$ws.Range("I126").Value = 932.5
$ws.Range("K126").Value = 2797.5
$ws.Range("M126").Value = 2142.5
$ws.Range("H131").Value = 868.99
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 868.99
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2606.97
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12686.97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2116.5833
$ws.Range("I97").Value = 2140.9
$ws.Range("J97").Value = 1995
$ws.Range("K97").Value = 2140.9
$ws.Range("L97").Value = 1995
$ws.Range("M97").Value = -1644.9
$ws.Range("N97").Value = -2987
$ws.Range("H122").Value = 2036.0526
$ws.Range("I122").Value = 1999
$ws.Range("K122").Value = 5997
$ws.Range("M122").Value = -3547
$ws.Range("H126").Value = 11518.148
$ws.Range("I126").Value = 2799.3333
$ws.Range("J126").Value = 22416.666
$ws.Range("K126").Value = 8397.999899999999
$ws.Range("L126").Value = 67249.99800000001
$ws.Range("M126").Value = -5927.999899999999
$ws.Range("N126").Value = -72189.99800000001
$ws.Range("H132").Value = 4367.7144
$ws.Range("I132").Value = 4013.6
$ws.Range("J132").Value = 5253
$ws.Range("K132").Value = 12040.8
$ws.Range("L132").Value = 15759
$ws.Range("M132").Value = -9510.799999999999
$ws.Range("N132").Value = -20819
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 34749.668
$ws.Range("J62").Value = 34749.668
$ws.Range("L62").Value = 34749.668
$ws.Range("N62").Value = -35997.668
$ws.Range("H64").Value = 30929.6
$ws.Range("J64").Value = 30929.6
$ws.Range("L64").Value = 30929.6
$ws.Range("N64").Value = -31379.6
$ws.Range("H65").Value = 34749.668
$ws.Range("J65").Value = 34749.668
$ws.Range("L65").Value = 104249.004
$ws.Range("N65").Value = -110489.004
$ws.Range("H67").Value = 30929.6
$ws.Range("J67").Value = 30929.6
$ws.Range("L67").Value = 30929.6
$ws.Range("N67").Value = -32489.6
